# Update TPM-derived NATMI ligand-receptor metrics (Thbs1-Cd47) with
# recomputed values from the new TPM-based script run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("M2").Value = 32.51511900000001
$ws.Range("N2").Value = 97.54535700000001
$ws.Range("O2").Value = 0.218203973858649
$ws.Range("P2").Value = 0.2182039738586489
$ws.Range("Q2").Value = 88.68447766896902
$ws.Range("R2").Value = 798.1602990207211
$ws.Range("S2").Value = 0.006755702940760345
$ws.Range("T2").Value = 0.006755702940760345

# Row 3
$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("O3").Value = 0.6017421411306194
$ws.Range("P3").Value = 0.6017421411306194
$ws.Range("Q3").Value = 244.565607737949
$ws.Range("R3").Value = 2201.090469641541
$ws.Range("S3").Value = 0.01863023427359281
$ws.Range("T3").Value = 0.01863023427359282

# Row 4
$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 26.830279
$ws.Range("N4").Value = 80.490837
$ws.Range("O4").Value = 0.1800538850107317
$ws.Range("P4").Value = 0.1800538850107317
$ws.Range("Q4").Value = 73.17916563146234
$ws.Range("R4").Value = 658.612490683161
$ws.Range("S4").Value = 0.005574557323370723
$ws.Range("T4").Value = 0.005574557323370722

# Row 5
$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("M5").Value = 32.51511900000001
$ws.Range("N5").Value = 97.54535700000001
$ws.Range("O5").Value = 0.218203973858649
$ws.Range("P5").Value = 0.2182039738586489
$ws.Range("Q5").Value = 1714.656338897502
$ws.Range("R5").Value = 15431.90705007752
$ws.Range("S5").Value = 0.1306170952973476
$ws.Range("T5").Value = 0.1306170952973476

# Row 6
$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("O6").Value = 0.6017421411306194
$ws.Range("P6").Value = 0.6017421411306194
$ws.Range("Q6").Value = 4728.515977164342
$ws.Range("R6").Value = 42556.64379447908
$ws.Range("S6").Value = 0.3602033876954195
$ws.Range("T6").Value = 0.3602033876954195

# Row 7
$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 26.830279
$ws.Range("N7").Value = 80.490837
$ws.Range("O7").Value = 0.1800538850107317
$ws.Range("P7").Value = 0.1800538850107317
$ws.Range("Q7").Value = 1414.871277627449
$ws.Range("R7").Value = 12733.84149864704
$ws.Range("S7").Value = 0.1077804177495836
$ws.Range("T7").Value = 0.1077804177495836

# Row 8
$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("M8").Value = 32.51511900000001
$ws.Range("N8").Value = 97.54535700000001
$ws.Range("O8").Value = 0.218203973858649
$ws.Range("P8").Value = 0.2182039738586489
$ws.Range("Q8").Value = 1061.099141293737
$ws.Range("R8").Value = 9549.892271643635
$ws.Range("S8").Value = 0.08083117562054097
$ws.Range("T8").Value = 0.08083117562054096

# Row 9
$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("O9").Value = 0.6017421411306194
$ws.Range("P9").Value = 0.6017421411306194
$ws.Range("Q9").Value = 2926.198171109278
$ws.Range("R9").Value = 26335.7835399835
$ws.Range("S9").Value = 0.222908519161607
$ws.Range("T9").Value = 0.222908519161607

# Row 10
$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 26.830279
$ws.Range("N10").Value = 80.490837
$ws.Range("O10").Value = 0.1800538850107317
$ws.Range("P10").Value = 0.1800538850107317
$ws.Range("Q10").Value = 875.5799419824172
$ws.Range("R10").Value = 7880.219477841754
$ws.Range("S10").Value = 0.06669890993777731
$ws.Range("T10").Value = 0.0666989099377773

